$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "2024-09-04 17:18:27"
# B2:G2 are present but empty (mirrors the source export which always emits
# all 9 columns, even when blank). Leading "'" keeps them literal empty text
# instead of Excel collapsing them back down to "no cell".
"B2","C2","D2","E2","F2","G2" | ForEach-Object { $ws.Range($_).Value = "'" }
$ws.Range("H2").Value = 238001
$ws.Range("I2").Value = "kit roletes"

# --- Row 3 ---
$ws.Range("A3").Value = "2024-09-04 17:24:19"
"B3","C3","D3","E3","F3","G3" | ForEach-Object { $ws.Range($_).Value = "'" }
$ws.Range("H3").Value = 238001
$ws.Range("I3").Value = "kit roletes"

# --- Row 4 ---
# A4 and H4 must stay plain text (not auto-converted to a date / number).
$ws.Range("A4").Value = "'2024-09-04"
$ws.Range("B4").Value = "17:28:07"
$ws.Range("C4").Value = "X3BK036141"
$ws.Range("D4").Value = "WF-M5799"
$ws.Range("E4").Value = "EPSON"
$ws.Range("F4").Value = "BL. A (3º) - UTI A PST 3"
$ws.Range("G4").Value = "10.16.13.143"
$ws.Range("H4").Value = "'238001"
$ws.Range("I4").Value = "kit roletes"
